$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update direction (column C) for PB4 row: DO -> PWM
$ws.Range("C4").Value = "PWM"

# Update description (column F) for PB4 row: add dutycycle range info
$ws.Range("F4").Value = "使能RF，可脉冲供电1~100%"

# Widen column F to fit new text
$ws.Columns("F").ColumnWidth = 24.53125

# Update selection to reflect last edited cell
$ws.Range("E14").Select()
